$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.677.12"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "1.798.97"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").Value = "313.79"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").Value = "0.5363"
$ws.Range("E7").Value = "  -1.47%  "
$ws.Range("D8").Value = "0.3779"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").Value = "0.07539"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").Value = "42.47"
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("D11").Value = "1.117"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "20.99"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "6.184"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "7.431"
$ws.Range("E15").Value = "  +4.05%  "
$ws.Range("D16").Value = "1.795.07"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "90.31"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "0.00001066"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").Value = "0.06452"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").Value = "17.25"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").Value = "5.940"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("D23").Value = "28.660.45"
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("D24").Value = "11.21"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "2.097"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "160.58"
$ws.Range("E26").Value = "  +2.85%  "
$ws.Range("D27").Value = "20.45"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").Value = "2.379"
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("D29").Value = "2.001.67"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "122.99"
$ws.Range("D31").Value = "1.108"
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("D32").Value = "0.1027"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").Value = "5.669"
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("D34").Value = "3.693"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("D35").Value = "0.2256"
$ws.Range("E35").Value = "  +6.76%  "
$ws.Range("D36").Value = "0.06478"
$ws.Range("E36").Value = "  +6.86%  "
$ws.Range("D37").Value = "8.957"
$ws.Range("E37").Value = "  +3.64%  "
$ws.Range("D38").Value = "0.02306"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").Value = "5.044"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").Value = "11.37"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").Value = "1.211"
$ws.Range("E41").Value = "  +5.14%  "
$ws.Range("D42").Value = "0.6261"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").Value = "1.394"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").Value = "13.48"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").Value = "0.5889"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").Value = "3.661"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "126.40"
$ws.Range("E48").Value = "  +3.33%  "
$ws.Range("D49").Value = "1.965"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("E51").Value = "  +1.77%  "
